# Generate Report for Handoff
# Updates the localization-status workbook so that the first entry (row 2)
# moves from "handed back" to "ready for handoff", and the second entry's
# (row 3) GUID file name changes. The per-language sheets lose their
# "Latest Target File" / "Latest Handback File" columns (F/G) because the
# file has not been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "21087872-3067-4e97-a60a-60def8703ad0"
$oldGuid2 = "d3e5deb1-5f73-4b6b-844d-b4a9e5f12179"
$newGuid1 = "a3573822-9d94-4592-834d-ffb5b55a6ca7"
$newGuid2 = "ffffb6ad0d8b-c1c0-4c1f-bfc6-69772c211997"
$newHash  = "2924ee7c1e8c42dfa5b48b47664a134ed1cea41e"

$status = "Ready for handoff"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2,1).Value = ($newGuid1 + ".md")
$wsOverview.Cells.Item(2,2).Value = $status
$wsOverview.Cells.Item(2,3).Value = $status
$wsOverview.Cells.Item(2,4).Value = "2016-02-13 05:02:08"

$wsOverview.Cells.Item(3,1).Value = ($newGuid2 + ".md")
$wsOverview.Cells.Item(3,2).Value = $status
$wsOverview.Cells.Item(3,3).Value = $status
$wsOverview.Cells.Item(3,4).Value = "2016-02-13 05:02:08"

# rebuild hyperlinks (clear whole-sheet collection, then re-add in order)
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(2,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid1 + ".md"), "", "", ($newGuid1 + ".md"))
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid2 + ".md"), "", "", ($newGuid2 + ".md"))

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(2,1).Value = ($newGuid1 + ".md")
$wsZhCn.Cells.Item(2,2).Value = ".md"
$wsZhCn.Cells.Item(2,3).Value = $status
$wsZhCn.Cells.Item(2,4).Value = ($newGuid1 + "." + $newHash + ".zh-cn.xlf")
$wsZhCn.Cells.Item(2,5).Value = "2016-03-13 05:02:04"
$wsZhCn.Cells.Item(2,6).Clear()
$wsZhCn.Cells.Item(2,7).Clear()
$wsZhCn.Cells.Item(2,8).Value = "0001-01-01 00:00:00"

$wsZhCn.Cells.Item(3,1).Value = ($newGuid2 + ".md")
$wsZhCn.Cells.Item(3,2).Value = ".md"
$wsZhCn.Cells.Item(3,3).Value = $status
$wsZhCn.Cells.Item(3,4).Value = ($newGuid1 + "." + $newHash + ".zh-cn.xlf")
$wsZhCn.Cells.Item(3,5).Value = "2016-03-13 05:02:04"
$wsZhCn.Cells.Item(3,6).Clear()
$wsZhCn.Cells.Item(3,7).Clear()
$wsZhCn.Cells.Item(3,8).Value = "0001-01-01 00:00:00"

$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid1 + ".md"), "", "", ($newGuid1 + ".md"))
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2,2), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid1 + ".md"), "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2,4), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/974d58029c702ab535d91f93d2ed56cff54afb84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newGuid1 + "." + $newHash + ".zh-cn.xlf"), "", "", ($newGuid1 + "." + $newHash + ".zh-cn.xlf"))
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid2 + ".md"), "", "", ($newGuid2 + ".md"))
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3,2), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid2 + ".md"), "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3,4), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/974d58029c702ab535d91f93d2ed56cff54afb84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $newGuid1 + "." + $newHash + ".zh-cn.xlf"), "", "", ($newGuid1 + "." + $newHash + ".zh-cn.xlf"))

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(2,1).Value = ($newGuid1 + ".md")
$wsDeDe.Cells.Item(2,2).Value = ".md"
$wsDeDe.Cells.Item(2,3).Value = $status
$wsDeDe.Cells.Item(2,4).Value = ($newGuid1 + "." + $newHash + ".de-de.xlf")
$wsDeDe.Cells.Item(2,5).Value = "2016-03-13 05:02:08"
$wsDeDe.Cells.Item(2,6).Clear()
$wsDeDe.Cells.Item(2,7).Clear()
$wsDeDe.Cells.Item(2,8).Value = "0001-01-01 00:00:00"

$wsDeDe.Cells.Item(3,1).Value = ($newGuid2 + ".md")
$wsDeDe.Cells.Item(3,2).Value = ".md"
$wsDeDe.Cells.Item(3,3).Value = $status
$wsDeDe.Cells.Item(3,4).Value = ($newGuid1 + "." + $newHash + ".de-de.xlf")
$wsDeDe.Cells.Item(3,5).Value = "2016-03-13 05:02:08"
$wsDeDe.Cells.Item(3,6).Clear()
$wsDeDe.Cells.Item(3,7).Clear()
$wsDeDe.Cells.Item(3,8).Value = "0001-01-01 00:00:00"

$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid1 + ".md"), "", "", ($newGuid1 + ".md"))
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2,2), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid1 + ".md"), "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2,4), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6cfa1ddf589e6997c5387ace23554f9178a56773/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newGuid1 + "." + $newHash + ".de-de.xlf"), "", "", ($newGuid1 + "." + $newHash + ".de-de.xlf"))
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3,1), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid2 + ".md"), "", "", ($newGuid2 + ".md"))
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3,2), ("https://github.com/OpenLocalizationTest/oltest/blob/8a5bf05bb058937f55e0d2f6f912cf05e8af4110/e2e/" + $newGuid2 + ".md"), "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3,4), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6cfa1ddf589e6997c5387ace23554f9178a56773/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $newGuid1 + "." + $newHash + ".de-de.xlf"), "", "", ($newGuid1 + "." + $newHash + ".de-de.xlf"))

$wb.Save()
